# Append a new partnership-email record row to the report and rename the
# worksheet to reflect its content, mirroring the UiPath
# AppendMetadataToExcel.xaml workflow output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the default "Sheet1" tab.
$ws.Name = "Partnership_Emails"

# --- New data row (row 2) ------------------------------------------------
# A: Date Processed (date/time value, formatted like a short date+time).
# Using the raw serial number (instead of a [DateTime] object) avoids Excel
# auto-applying its own default datetime format before ours is set, and the
# "m/d/yy h:mm" code maps onto Excel's built-in numFmtId 22 instead of
# registering a redundant custom format.
$ws.Range("A2").Value = 45854.104351851849
$ws.Range("A2").NumberFormat = "m/d/yy h:mm"

# B: Sender Name
$ws.Range("B2").Value = "Moris Mwai"

# C: Company
$ws.Range("C2").Value = "Tech-Neo GmbH"

# D: Address
$ws.Range("D2").Value = "Am main City, Germany"

# E: VAT ID (multi-line text copied from the source email, wrapped)
$ws.Range("E2").Value = "DE1567890`n`n`n`n"
$ws.Range("E2").WrapText = $true

# F: Email
$ws.Range("F2").Value = "morismwai1@gmail.com"

# G: Email Subject
$ws.Range("G2").Value = "Partnership Offer"

# Row 2 is tall enough to show the wrapped VAT ID text.
$ws.Rows.Item(2).RowHeight = 60

# --- Column sizing ---------------------------------------------------------
# Auto-fit the newly populated columns (C:G) to their content, same as the
# existing bestFit columns A & B already in the workbook.
$ws.Range("C1:G2").EntireColumn.AutoFit()

# --- Selection ---------------------------------------------------------
# After the append, the active cell moves to the first empty cell past the
# new data.
$ws.Range("H1").Select()
